$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 65569
$ws.Range("B2").Value = "Pedro Lucas Almeida"
$ws.Range("C2").Value = "Marketing"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 45095
$ws.Range("G2").Value = 11025.98

# Row 3
$ws.Range("A3").Value = 47595
$ws.Range("B3").Value = "Luiza Fogaça"
$ws.Range("C3").Value = "TI"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45083
$ws.Range("G3").Value = 5941.98

# Row 4
$ws.Range("A4").Value = 6630
$ws.Range("B4").Value = "Dr. Augusto Cardoso"
$ws.Range("C4").Value = "Recursos Humanos"
$ws.Range("D4").Value = "Doença"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 45099
$ws.Range("G4").Value = 11723.4

# Row 5
$ws.Range("A5").Value = 66020
$ws.Range("B5").Value = "Ana Carolina Cunha"
$ws.Range("C5").Value = "TI"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45089
$ws.Range("G5").Value = 11596.75

# Row 6
$ws.Range("A6").Value = 62288
$ws.Range("B6").Value = "Gustavo Melo"
$ws.Range("C6").Value = "Vendas"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 45082
$ws.Range("G6").Value = 5633.98

# Row 7
$ws.Range("A7").Value = 64238
$ws.Range("B7").Value = "Rebeca Mendes"
$ws.Range("C7").Value = "Operações"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45084
$ws.Range("G7").Value = 12048.72

# Row 8
$ws.Range("A8").Value = 95674
$ws.Range("B8").Value = "Vinicius Melo"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45096
$ws.Range("G8").Value = 10461.03

# Row 9
$ws.Range("A9").Value = 2844
$ws.Range("B9").Value = "Eloah Santos"
$ws.Range("C9").Value = "Marketing"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45095
$ws.Range("G9").Value = 10638

# Row 10
$ws.Range("A10").Value = 75037
$ws.Range("B10").Value = "Cauã Melo"
$ws.Range("C10").Value = "Jurídico"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45083
$ws.Range("G10").Value = 3386.9

# Row 11
$ws.Range("A11").Value = 90143
$ws.Range("B11").Value = "João Lucas Barros"
$ws.Range("C11").Value = "P&D"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45083
$ws.Range("G11").Value = 8105.05
